$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G9").Value = 176
$ws.Range("G12").Value = 176
$ws.Range("G14").Value = 179
$ws.Range("G16").Value = 27
$ws.Range("G17").Value = 31
$ws.Range("G18").Value = 118
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
